$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.452.00'
$ws.Range("E2").Value = '  +0.54%  '
$ws.Range("D3").Value = '3.345.23'
$ws.Range("E3").Value = '  +0.76%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '186.51'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.07%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("E8").Value = '  -1.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.183'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.21%  '
$ws.Range("E10").Value = '  -0.21%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '47.11'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.01%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000270'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.49%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '648.44'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +7.02%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '3.873.48'
$ws.Range("E14").Value = '  +0.62%  '
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.49'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.99%  '
$ws.Range("D16").Value = '66.462.09'
$ws.Range("E16").Value = '  +0.50%  '
$ws.Range("E17").Value = '  +0.13%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.91'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.31%  '
$ws.Range("D19").Value = '3.339.96'
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.16'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.70%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.901'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.45%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.75'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.09%  '
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '100.23'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.52%  '
$ws.Range("E25").Value = '  +0.82%  '
$ws.Range("E26").Value = '  +1.45%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.57'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '32.18'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.11%  '
$ws.Range("E29").Value = '  -2.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.88'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.43%  '
$ws.Range("B31").Value = 'dogwifhat'
$ws.Range("C31").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.99'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.04%  '
$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '603.86'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.12'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.36%  '
$ws.Range("D34").Value = '3.879.09'
$ws.Range("E34").Value = '  +4.59%  '
$ws.Range("E35").Value = '  +0.91%  '
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '56.67'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.96%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.73'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.61%  '
$ws.Range("E39").Value = '  -0.62%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '33.39'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.09%  '
$ws.Range("E41").Value = '  -1.83%  '
$ws.Range("E42").Value = '  -2.08%  '
$ws.Range("E43").Value = '  +0.07%  '
$ws.Range("E44").Value = '  -1.50%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0418'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.41%  '
$ws.Range("E46").Value = '  -0.91%  '
$ws.Range("E47").Value = '  +0.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.55'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.05%  '
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.37'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +9.87%  '
$ws.Range("B50").Value = 'CoreDAO'
$ws.Range("C50").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.85'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -15.16%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '129.55'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.14%  '
